$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").EntireRow.Insert()

$ws.Range("D11").NumberFormat = $ws.Range("D12").NumberFormat

$ws.Range("A11").Value = 12
$ws.Range("B11").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = 44483
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = 100112028
$ws.Range("G11").Value = "Sandia"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 800
$ws.Range("L11").Value = 800
$ws.Range("M11").Value = 800
$ws.Range("N11").Value = "$/kilo (volumen en unidades)"
$ws.Range("O11").Value = "Perú"
$ws.Range("P11").Value = 800
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = "Hortaliza"
